$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 137, shifting existing rows (and the rest of the
# data block) down by 4. This turns the old A1:T168 range into A1:T172.
$ws.Range("A137:T140").Insert()

# Fill in the new block (previously-missing week, date 2022-03-24 / serial 44644)
# Column layout:
# A Mercado ID | B Mercado | C Region | D Fecha | E Codreg | F Tipo
# G Producto ID | H Producto | I Categoria ID | J Categoria | K Variedad
# L Calidad | M Volumen | N Precio minimo | O Precio maximo
# P Precio promedio ponderado | Q Unidad de comercializacion | R Origen
# S Precio $/Kg | T Kg / unidad

$mercado = "Agr" + [char]0x00ED + "cola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$tipo = "Fruta"
$producto = "Tropicales y subtropicales"
$categoria = "Pi" + [char]0x00F1 + "a"
$variedad = "Caramelo"
$origen = "Ecuador"

$rows = @(
    @{ row = 137; calidad = "Especial"; vol = 200; pmin = 19000; pmax = 20000; pprom = 19500; unidad = "`$/caja 10 unidades"; precioKg = 1950; kgUnidad = 10 },
    @{ row = 138; calidad = "Primera";  vol = 260; pmin = 19000; pmax = 20000; pprom = 19500; unidad = "`$/caja 12 unidades"; precioKg = 1625; kgUnidad = 12 },
    @{ row = 139; calidad = "Segunda";  vol = 270; pmin = 19000; pmax = 20000; pprom = 19500; unidad = "`$/caja 14 unidades"; precioKg = 1393; kgUnidad = 14 },
    @{ row = 140; calidad = "Tercera";  vol = 270; pmin = 19000; pmax = 20000; pprom = 19500; unidad = "`$/caja 16 unidades"; precioKg = 1219; kgUnidad = 16 }
)

foreach ($r in $rows) {
    $row = $r.row
    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = 44644
    $ws.Cells.Item($row, 5).Value = 15
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = 100108
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = 100108005
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.calidad
    $ws.Cells.Item($row, 13).Value = $r.vol
    $ws.Cells.Item($row, 14).Value = $r.pmin
    $ws.Cells.Item($row, 15).Value = $r.pmax
    $ws.Cells.Item($row, 16).Value = $r.pprom
    $ws.Cells.Item($row, 17).Value = $r.unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.precioKg
    $ws.Cells.Item($row, 20).Value = $r.kgUnidad
}
